# Tee valitut raportit vain Järjestelmävastaavan saataville
# (Make the selected reports available only to the System Administrator)
#
# For the listed report rows, clear every role column except
# "Järjestelmävastaava" (column D), leaving a single blank-space marker
# in column G (matching the convention already used by untouched rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oikeudet")

$rows = @(68, 73, 80, 82, 91)

foreach ($r in $rows) {
    # Clear every role column from E through X (column D - Järjestelmävastaava - stays untouched)
    $ws.Range("E$r`:X$r").Value = ""
    # Column G keeps a single-space placeholder, same convention as the other rows
    $ws.Range("G$r").Value = " "
}

# Leave the selection where the edit happened, matching the authored view state
$ws.Activate()
$ws.Range("D82:X82").Select()
